$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teste")
$ws.Range("A2").Value = "BarbeiroTeste"
$ws.Range("A3").Value = "BarbeiroTeste"
$ws.Range("E3").Value = "setRg"
$ws.Range("F3").Value = "setRg"
